$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I ("Polygons" is H, "Generations" was I).
# This shifts old I:M -> J:N and makes room for the new "HOF" column.
$ws.Columns("I:I").Insert()

# New column header
$ws.Range("I1").Value = "HOF"

# New column values for the existing 8 data rows (rows 2-9)
$ws.Range("I2:I9").Value = 20

# Update shape text for existing rows: "Regular Polygon" -> "Regular Polygon (5)"
$ws.Range("F2:F9").Value = "Regular Polygon (5)"

# Fill in the new 9th data row (row 10) with its values
$ws.Range("B10").Value = "Mona Lisa (Color)"
$ws.Range("C10").Value = "Tournment (size 2)"
$ws.Range("D10").Value = "SimulatedBinaryBounded (0.9)"
$ws.Range("E10").Value = "PolynomialBounded (0.01)"
$ws.Range("F10").Value = "Regular Polygon (6)"
$ws.Range("G10").Value = 800
$ws.Range("H10").Value = 700
$ws.Range("I10").Value = 20
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 1533.71795852269
$ws.Range("L10").Value = "31841.32s"
$ws.Range("M10").Value = 0.59226276708125103
$ws.Range("N10").Value = "35567.84s"

# Update the visible selection to match the authored workbook
$ws.Range("B11:B12").Select()
